# Apply the grade-sheet edits to "Лист1":
#  - Fill in previously-blank lab scores (column D / column E) with 0
#  - Correct Ибрамхалилов's "Лаба №1" score (column B, row 10) from 0 to 5
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Previously-empty "Лаба №2" (column D) cells -> 0
$ws.Range("D5").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("D16").Value = 0
$ws.Range("D18").Value = 0

# Previously-empty "Лаба №3" (column E) cell -> 0
$ws.Range("E15").Value = 0

# Ибрамхалилов Роман Ламетович, "Лаба №1" score corrected from 0 to 5
$ws.Range("B10").Value = 5
